# New Version of KWKG with variable full load hour calculation
$wb = $excel.ActiveWorkbook

# --- gen_economics: observation period 15 -> 20 years, interest rate 0.021 -> 0.02 ---
$gen = $wb.Worksheets.Item("gen_economics")
$gen.Range("B2").Value = 20
$gen.Range("B4").Value = 0.02

# --- gas_economics: widen price-change validity range placeholder 999 -> 9999 ---
$gas = $wb.Worksheets.Item("gas_economics")
$gas.Range("B2").Value = "[0 - 9999: 138]"
$gas.Range("C2").Value = "[0 - 9999: 0.058]"

# --- el_economics: same range widening for both tariff rows ---
$el = $wb.Worksheets.Item("el_economics")
$el.Range("B2").Value = "[0 - 9999: 73]"
$el.Range("C2").Value = "[0 - 9999: 0.258]"
$el.Range("B3").Value = "[0 - 9999: 73]"
$el.Range("C3").Value = "[0 - 9999: 0.1856]"

# --- pel_economics: same range widening ---
$pel = $wb.Worksheets.Item("pel_economics")
$pel.Range("B2").Value = "[0 - 9999: 1]"
$pel.Range("C2").Value = "[0 - 9999: 0.048]"

# --- further_parameters: drop the dT_max row (row 6) entirely ---
$further = $wb.Worksheets.Item("further_parameters")
$further.Rows.Item(6).Delete()

# --- reorder tabs: ep_table moves to sit right before further_parameters ---
$epTable = $wb.Worksheets.Item("ep_table")
$epTable.Move($further)

# --- restore per-sheet cell selections, and move the active tab to gen_economics ---
$dev = $wb.Worksheets.Item("dev_economics")
$comp = $wb.Worksheets.Item("comp_economics")

$gas.Range("B2").Select()
$el.Range("C3").Select()
$pel.Range("C9").Select()
$dev.Range("D21").Select()
$comp.Range("F20").Select()
$epTable.Range("N25").Select()
$further.Range("D11").Select()

$gen.Activate()
$gen.Range("D26").Select()
